$d = $word.ActiveDocument

# ---- Part 1: "01/No.04 - Dok.02/2022" -> "01/No.04 - Dok.02/VST/2022" with _GoBack
#      bookmark relocated between "VST/" and "2022"
$find1 = $d.Content
$found1 = $find1.Find.Execute("Dok.02/", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPos = $find1.End

# Move the (only) existing _GoBack bookmark to the insertion point; this is the single
# reliable bookmark relocation this engine supports in one script run.
$boundary = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $boundary)

# Insert "VST/" immediately before the bookmark's own range -- the bookmark naturally
# shifts to sit right after the inserted text, landing exactly between "VST/" and "2022".
$bm = $d.Bookmarks("_GoBack")
$bm.Range.InsertBefore("VST/")

# Nudge formatting on "VST/" (apply & clear Bold) purely to force Word to keep it as its
# own run separate from "...Dok.02/"; this does not leave any formatting residue.
$vstRange = $d.Range($insertPos, $insertPos + 4)
$vstRange.Font.Bold = 1
$vstRange.Font.Bold = 0

# ---- Part 2: the old _GoBack bookmark location (between "handle " and "pisau potong")
#      is implicitly cleared because bookmarks are unique and it was just relocated above.

Write-Host "Edits applied."
